$d = $word.ActiveDocument
$d.Content.Find.Execute("qr_code", $true, $true, $false, $false, $false,
                         $true, 1, $false, "QR_code", 2)
